# Insert a new weekly price record as row 21 (Macroferia Regional de Talca - Haba),
# pushing the existing rows 21-44 down to 22-45.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = 44494
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = 100112026
$ws.Range("G21").Value = "Haba"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 500
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 7000
$ws.Range("N21").Value = "`$/saco 25 kilos"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 280
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
